# "Add files via upload" — the metadata sheet's organization contact info
# (section 2, rows 7-10) was refreshed to a new contact person.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These are the unlocked "fill-in" cells (protection: locked=0), so they can
# be edited directly even though the worksheet itself is protected.
$ws.Range("B7").Value() = "Мааткулова Ж. Б."
$ws.Range("B8").Value() = "j.maatkulova@stat.kg"
$ws.Range("B9").Value() = " (312) 32 55 46"
$ws.Range("B10").Value() = "www.stat.gov.kg"

# Leave the cursor on B8, matching where the editor left off when saving.
$ws.Range("B8").Select() | Out-Null
